$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.758.13"
$ws.Range("E2").Value = "  +2.83%  "
$ws.Range("D3").Value = "'1.789.86"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'223.07"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'32.83"
$ws.Range("E8").Value = "  +6.88%  "
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("D10").Value = "'0.0689"
$ws.Range("E10").Value = "  +3.74%  "
$ws.Range("E11").Value = "  +1.49%  "
$ws.Range("D12").Value = "'2.046.84"
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("D13").Value = "'11.08"
$ws.Range("E13").Value = "  +10.72%  "
$ws.Range("D14").Value = "'1.793.81"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("D15").Value = "'34.736.07"
$ws.Range("E15").Value = "  +2.83%  "
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("E17").Value = "  +3.40%  "
$ws.Range("D18").Value = "'68.48"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").Value = "'253.62"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("D20").Value = "'0.0₃0786"
$ws.Range("E20").Value = "  +6.35%  "
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "'10.49"
$ws.Range("E22").Value = "  +2.11%  "
$ws.Range("E23").Value = "  +0.64%  "
$ws.Range("D24").Value = "'2.13"
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("D25").Value = "'158.98"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("E27").Value = "  +1.41%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("D31").Value = "'0.0515"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("D35").Value = "'1.438.05"
$ws.Range("E35").Value = "  -2.73%  "
$ws.Range("E36").Value = "  -1.10%  "
$ws.Range("E37").Value = "  +2.55%  "
$ws.Range("D38").Value = "'0.632"
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("D39").Value = "'82.98"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").Value = "'2.81"
$ws.Range("E40").Value = "  +4.20%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").Value = "'0.906"
$ws.Range("E42").Value = "  +2.19%  "
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "'0.0503"
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'1.06"
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("E46").Value = "  +4.35%  "
$ws.Range("D47").Value = "'1.943.32"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("D48").Value = "'105.04"
$ws.Range("E48").Value = "  +7.81%  "
$ws.Range("D49").Value = "'12.00"
$ws.Range("E49").Value = "  +1.06%  "
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").Value = "'49.81"
$ws.Range("E51").Value = "  -2.39%  "
